# Fixed up Myxicola and a few more loose ends in second review
#
# On the "Materials" sheet:
#   - Drop the Taxon_Local_ID column (A) entirely.
#   - Drop the suborder / infraorder / superfamily columns entirely.
# On the "Materials" sheet's "authority" mapping cell:
#   - ${summary.Author} -> ${summary.authority}
#
# (Deleting the now-unused header/value cells lets the shared-strings table
#  shrink on its own: suborder, infraorder, superfamily, ${iNaturalistTaxonId},
#  ${suborder}, ${infraorder}, ${superfamily} are all removed automatically
#  once nothing references them.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

# Fix the authority formula text before shifting columns around, while the
# layout still matches the original file (scientificNameAuthorship is BB).
# NOTE: single-quote the literal so PowerShell doesn't try to expand
# "${summary.authority}" as a variable reference.
$ws.Range("BB2").Value = '${summary.authority}'

# Remove the suborder / infraorder / superfamily columns (originally AR:AT).
$ws.Range("AR1:AT2").EntireColumn.Delete()

# Remove the Taxon_Local_ID column (A), shifting everything else left.
$ws.Columns.Item(1).Delete()
